# Update data citing and optimization function
#
# For both the "GUI" (Guinea) and "SL" (Sierra Leone) sheets:
#   - insert a new column D "Infectious" = Cases - Deaths (B - C),
#     pushing the existing "CFR (%)" column from D to E
#   - update selection / active-sheet state to match the edited workbook
#     (SL becomes the active tab)

$wb = $excel.ActiveWorkbook

# ---- GUI sheet -----------------------------------------------------------
$guiSheet = $wb.Worksheets.Item("GUI")
[void]$guiSheet.Activate()

# Insert a new column at D; existing D (CFR %) shifts to E.
$guiSheet.Columns.Item(4).Insert()

# Header for the new column.
$guiSheet.Range("D1").Value = "Infectious"

# Row 2 formula (kept distinct, matching how the workbook stores it),
# then fill the remaining rows as one formula fill so the saved file uses
# a shared formula group for D3:D18.
$guiSheet.Range("D2").Formula = "=B2-C2"
$guiSheet.Range("D3:D18").Formula = "=B3-C3"

# Leave column D selected on this sheet.
[void]$guiSheet.Range("D:D").Select()

# ---- SL sheet --------------------------------------------------------------
$slSheet = $wb.Worksheets.Item("SL")
[void]$slSheet.Activate()

$slSheet.Columns.Item(4).Insert()

$slSheet.Range("D1").Value = "Infectious"

$slSheet.Range("D2").Formula = "=B2-C2"
$slSheet.Range("D3:D18").Formula = "=B3-C3"

# Final selection/active sheet is SL, matching the committed workbook state.
[void]$slSheet.Range("F13").Select()
